# Update cryptocurrency price/volume data per the Tue Sep 10 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.491.76"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.366.96"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.26"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  +5.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.50"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").Value = "2.792.02"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "57.504.68"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "2.377.94"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.66"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "331.08"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.12%  "
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.74"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("B24").Value = "Kaspa"
$ws.Range("C24").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.166"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.69%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.70"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +13.21%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +12.54%  "
$ws.Range("D28").Value = "0.0₃0749"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.71"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.32"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.98%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.925"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.47%  "
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.62"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.71"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "151.36"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.388"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.68"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.37"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.05%  "
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0221"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.32"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +6.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.86"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.363"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.94%  "
